# Daily attendance processing - 2025-10-24 07:42:54
# Normalizes the "Recorded By" (column G) cell values: for a known set of
# multi-recorder strings, the last name in the comma-separated list is
# moved to the front (a right-rotation of the list), e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-value replacement map derived from the set of "Recorded By" values
# that need to be normalized. Only cells whose current value matches one
# of these keys (exactly) are touched; everything else is left alone.
$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
